$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 40; this shifts existing rows 40-79 down to 41-80
$ws.Rows.Item(40).Insert()

# Populate the new row 40 with the new data record
$ws.Range("A40").Value = 9
$ws.Range("B40").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C40").Value = "Metropolitana"
$ws.Range("D40").Value = 44669
$ws.Range("E40").Value = 13
$ws.Range("F40").Value = 100114002
$ws.Range("G40").Value = "Camote"
$ws.Range("H40").Value = "Sin especificar"
$ws.Range("I40").Value = "Primera"
$ws.Range("J40").Value = 1420
$ws.Range("K40").Value = 10000
$ws.Range("L40").Value = 11000
$ws.Range("M40").Value = 10500
$ws.Range("N40").Value = "$/malla 18 kilos"
$ws.Range("O40").Value = "Perú"
$ws.Range("P40").Value = 583
$ws.Range("Q40").Value = 18
$ws.Range("R40").Value = "Hortaliza"
